{"js": "// Deny_Privileges_Template.docx edit\n//\n// The only substantive content change in the target revision is inside the\n// \"petition_incomplete\" reason paragraph: the phrase\n//   \"provide the further information\"\n// is replaced with\n//   \"resubmit a complete and/or legible petition\"\n// (everything else in the upstream diff is just Word re-running its\n// grammar checker over already-unchanged text, which only inserts/repositions\n// <w:proofErr> bookkeeping elements and splits runs around them without any\n// visible text difference).\n\nconst body = context.document.body;\n\nconst oldText = \"provide the further information\";\nconst newText = \"resubmit a complete and/or legible petition\";\n\nconst results = body.search(oldText, { matchCase: true });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(`Could not find target text to replace: \"${oldText}\"`);\n}\n\n// Replace in place (there is exactly one occurrence in this template).\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Deny_Privileges_Template.docx edit\n#\n# The only substantive content change in the target revision is inside the\n# \"petition_incomplete\" reason paragraph: the phrase\n#   \"provide the further information\"\n# is replaced with\n#   \"resubmit a complete and/or legible petition\"\n# (everything else in the upstream diff is just Word re-running its\n# grammar checker over already-unchanged text, which only inserts/repositions\n# <w:proofErr> bookkeeping elements and splits runs around them without any\n# visible text difference).\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"provide the further information\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"resubmit a complete and/or legible petition\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute([ref]$find.Text, [ref]$find.MatchCase, [ref]$find.MatchWholeWord, [ref]$find.MatchWildcards, $null, [ref]$find.Forward, [ref]$find.Wrap, $null, $null, [ref]$find.Replacement.Text, 2)\n"}
